# Update "想去人数" (want-to-go count) values for the 苏州-漫展信息 workbook.
# Both the "展览" and "全部类型" sheets contain the same rows, so the same
# three cells (F12, F13, F26) are bumped by 1 on each of those sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F12").Value = 13697
    $ws.Range("F13").Value = 14146
    $ws.Range("F26").Value = 5247
}
